$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.12165933333333
$ws.Range("H2").Value = 45.364978
$ws.Range("I2").Value = 0.1696222886509932
$ws.Range("J2").Value = 0.1696222886509932
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 232.2511563162795
$ws.Range("R2").Value = 2090.260406846516
$ws.Range("S2").Value = 0.0171665018877096
$ws.Range("T2").Value = 0.0171665018877096
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.12165933333333
$ws.Range("H3").Value = 45.364978
$ws.Range("I3").Value = 0.1696222886509932
$ws.Range("J3").Value = 0.1696222886509932
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 765.0435276889702
$ws.Range("R3").Value = 6885.391749200732
$ws.Range("S3").Value = 0.05654706469735737
$ws.Range("T3").Value = 0.05654706469735736
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.12165933333333
$ws.Range("H4").Value = 45.364978
$ws.Range("I4").Value = 0.1696222886509932
$ws.Range("J4").Value = 0.1696222886509932
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 913.002794628629
$ws.Range("R4").Value = 8217.025151657661
$ws.Range("S4").Value = 0.06748325582557241
$ws.Range("T4").Value = 0.06748325582557239
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.12165933333333
$ws.Range("H5").Value = 45.364978
$ws.Range("I5").Value = 0.1696222886509932
$ws.Range("J5").Value = 0.1696222886509932
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.29665199999999
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 384.5773266059617
$ws.Range("R5").Value = 3461.195939453656
$ws.Range("S5").Value = 0.02842546624035387
$ws.Range("T5").Value = 0.02842546624035387
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 30.93224
$ws.Range("H6").Value = 92.79671999999999
$ws.Range("I6").Value = 0.3469723279862584
$ws.Range("J6").Value = 0.3469723279862583
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 475.0833456230932
$ws.Range("R6").Value = 4275.75011060784
$ws.Range("S6").Value = 0.03511508523278153
$ws.Range("T6").Value = 0.03511508523278152
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 30.93224
$ws.Range("H7").Value = 92.79671999999999
$ws.Range("I7").Value = 0.3469723279862584
$ws.Range("J7").Value = 0.3469723279862583
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 1564.941352484853
$ws.Range("R7").Value = 14084.47217236368
$ws.Range("S7").Value = 0.1156703333911582
$ws.Range("T7").Value = 0.1156703333911581
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 30.93224
$ws.Range("H8").Value = 92.79671999999999
$ws.Range("I8").Value = 0.3469723279862584
$ws.Range("J8").Value = 0.3469723279862583
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 1867.600700530933
$ws.Range("R8").Value = 16808.4063047784
$ws.Range("S8").Value = 0.1380409529909617
$ws.Range("T8").Value = 0.1380409529909616
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 30.93224
$ws.Range("H9").Value = 92.79671999999999
$ws.Range("I9").Value = 0.3469723279862584
$ws.Range("J9").Value = 0.3469723279862583
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.29665199999999
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 786.6754502868265
$ws.Range("R9").Value = 7080.079052581439
$ws.Range("S9").Value = 0.058145956371357
$ws.Range("T9").Value = 0.058145956371357
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.447555666666666
$ws.Range("H10").Value = 13.342667
$ws.Range("I10").Value = 0.04988900718188559
$ws.Range("J10").Value = 0.04988900718188558
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 68.30929884046377
$ws.Range("R10").Value = 614.7836895641739
$ws.Range("S10").Value = 0.005048981137885276
$ws.Range("T10").Value = 0.005048981137885276
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.447555666666666
$ws.Range("H11").Value = 13.342667
$ws.Range("I11").Value = 0.04988900718188559
$ws.Range("J11").Value = 0.04988900718188558
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 225.0132476744331
$ws.Range("R11").Value = 2025.119229069898
$ws.Range("S11").Value = 0.01663152253891306
$ws.Range("T11").Value = 0.01663152253891306
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.447555666666666
$ws.Range("H12").Value = 13.342667
$ws.Range("I12").Value = 0.04988900718188559
$ws.Range("J12").Value = 0.04988900718188558
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 268.5307652700544
$ws.Range("R12").Value = 2416.77688743049
$ws.Range("S12").Value = 0.0198480557084459
$ws.Range("T12").Value = 0.0198480557084459
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.447555666666666
$ws.Range("H13").Value = 13.342667
$ws.Range("I13").Value = 0.04988900718188559
$ws.Range("J13").Value = 0.04988900718188558
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.29665199999999
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 113.1112023167649
$ws.Range("R13").Value = 1018.000820850884
$ws.Range("S13").Value = 0.008360447796641355
$ws.Range("T13").Value = 0.008360447796641355
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 38.64755633333333
$ws.Range("H14").Value = 115.942669
$ws.Range("I14").Value = 0.4335163761808628
$ws.Range("J14").Value = 0.4335163761808628
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 15.35884066666667
$ws.Range("N14").Value = 46.076522
$ws.Range("O14").Value = 0.1012042817263867
$ws.Range("P14").Value = 0.1012042817263867
$ws.Range("Q14").Value = 593.5816598796908
$ws.Range("R14").Value = 5342.234938917218
$ws.Range("S14").Value = 0.04387371346801026
$ws.Range("T14").Value = 0.04387371346801026
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 38.64755633333333
$ws.Range("H15").Value = 115.942669
$ws.Range("I15").Value = 0.4335163761808628
$ws.Range("J15").Value = 0.4335163761808628
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 50.59256466666667
$ws.Range("N15").Value = 151.777694
$ws.Range("O15").Value = 0.3333704853712116
$ws.Range("P15").Value = 0.3333704853712116
$ws.Range("Q15").Value = 1955.27899300281
$ws.Range("R15").Value = 17597.51093702529
$ws.Range("S15").Value = 0.144521564743783
$ws.Range("T15").Value = 0.1445215647437829
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 38.64755633333333
$ws.Range("H16").Value = 115.942669
$ws.Range("I16").Value = 0.4335163761808628
$ws.Range("J16").Value = 0.4335163761808628
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 60.37715666666667
$ws.Range("N16").Value = 181.13147
$ws.Range("O16").Value = 0.397844271305776
$ws.Range("P16").Value = 0.397844271305776
$ws.Range("Q16").Value = 2333.429563521493
$ws.Range("R16").Value = 21000.86607169343
$ws.Range("S16").Value = 0.172472006780796
$ws.Range("T16").Value = 0.172472006780796
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 38.64755633333333
$ws.Range("H17").Value = 115.942669
$ws.Range("I17").Value = 0.4335163761808628
$ws.Range("J17").Value = 0.4335163761808628
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.43221733333333
$ws.Range("N17").Value = 76.29665199999999
$ws.Range("O17").Value = 0.1675809615966257
$ws.Range("P17").Value = 0.1675809615966258
$ws.Range("Q17").Value = 982.8930520715763
$ws.Range("R17").Value = 8846.037468644188
$ws.Range("S17").Value = 0.07264909118827353
$ws.Range("T17").Value = 0.07264909118827353
